$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the dates in D2/D3 (set to 2022-10-12 = 44846) and D6/D7 (set to 2022-09-28 = 44832)
$ws.Range("D2").Value = 44846
$ws.Range("D3").Value = 44846
$ws.Range("D6").Value = 44832
$ws.Range("D7").Value = 44832
